$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)     # StartSceneConfig
$ws2 = $wb.Worksheets.Item(2)     # Robot区

# ---------------------------------------------------------------
# Sheet1 (StartSceneConfig): fill in the previously-empty OuterPort
# values for rows 9-11, then append four new rows (12-15).
# ---------------------------------------------------------------
$ws1.Range("H9").Value  = 10005
$ws1.Range("H10").Value = 10006
$ws1.Range("H11").Value = 10010

# Row 12: Router / Router1
$null = $ws1.Range("C11:H11").Copy()
$null = $ws1.Range("C12:H12").PasteSpecial(-4122)
$ws1.Range("C12").Value = 7
$ws1.Range("D12").Value = 3
$ws1.Range("E12").Value = 1
$ws1.Range("F12").Value = "Router"
$ws1.Range("G12").Value = "Router1"
$ws1.Range("H12").Value = 10007

# Row 13: Router / Router2
$null = $ws1.Range("C11:H11").Copy()
$null = $ws1.Range("C13:H13").PasteSpecial(-4122)
$ws1.Range("C13").Value = 8
$ws1.Range("D13").Value = 4
$ws1.Range("E13").Value = 1
$ws1.Range("F13").Value = "Router"
$ws1.Range("G13").Value = "Router2"
$ws1.Range("H13").Value = 10008

# Row 14: Router / Router3
$null = $ws1.Range("C11:H11").Copy()
$null = $ws1.Range("C14:H14").PasteSpecial(-4122)
$ws1.Range("C14").Value = 9
$ws1.Range("D14").Value = 5
$ws1.Range("E14").Value = 1
$ws1.Range("F14").Value = "Router"
$ws1.Range("G14").Value = "Router3"
$ws1.Range("H14").Value = 10009

# Row 15: Map / Map1_1
$null = $ws1.Range("C11:H11").Copy()
$null = $ws1.Range("C15:H15").PasteSpecial(-4122)
$ws1.Range("C15").Value = 10
$ws1.Range("D15").Value = 1
$ws1.Range("E15").Value = 1
$ws1.Range("F15").Value = "Map"
$ws1.Range("G15").Value = "Map1_1"
$ws1.Range("H15").Value = 10011

# ---------------------------------------------------------------
# Sheet2 (Robot区): Zone (D6) 1 -> 2
# ---------------------------------------------------------------
$ws2.Range("D6").Value = 2

# ---------------------------------------------------------------
# Final selection / active-sheet state:
#   sheet2 row 6 fully selected (no longer the active tab)
#   sheet1 becomes the active tab, with H18 selected
# ---------------------------------------------------------------
$null = $ws2.Rows.Item(6).Select()
$null = $ws1.Range("H18").Select()
